$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "64.728.91"
$ws.Range("E2").Value = "  -2.22%  "

$ws.Range("D3").Value = "3.379.66"
$ws.Range("E3").Value = "  -3.93%  "

$ws.Range("E4").Value = "  +0.02%  "

Set-TextValue $ws.Range("D5") "557.91"
$ws.Range("E5").Value = "  -4.33%  "

Set-TextValue $ws.Range("D6") "177.00"
$ws.Range("E6").Value = "  -0.88%  "

Set-TextValue $ws.Range("D7") "0.618"
$ws.Range("E7").Value = "  -2.03%  "

$ws.Range("D8").Value = "3.372.82"
$ws.Range("E8").Value = "  -3.92%  "

$ws.Range("E9").Value = "  -0.17%  "

Set-TextValue $ws.Range("D10") "0.630"
$ws.Range("E10").Value = "  -1.46%  "

Set-TextValue $ws.Range("D11") "0.163"
$ws.Range("E11").Value = "  +0.00%  "

Set-TextValue $ws.Range("D12") "54.95"
$ws.Range("E12").Value = "  -1.57%  "

Set-TextValue $ws.Range("D13") "0.0000274"
$ws.Range("E13").Value = "  -2.02%  "

Set-TextValue $ws.Range("D14") "9.10"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "3.897.20"
$ws.Range("E15").Value = "  -4.58%  "

Set-TextValue $ws.Range("D16") "18.43"
$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("E17").Value = "  -2.55%  "

$ws.Range("D18").Value = "3.357.63"
$ws.Range("E18").Value = "  -5.00%  "

Set-TextValue $ws.Range("D19") "11.90"
$ws.Range("E19").Value = "  -1.63%  "

$ws.Range("D20").Value = "64.618.65"
$ws.Range("E20").Value = "  -2.42%  "

Set-TextValue $ws.Range("D21") "0.984"
$ws.Range("E21").Value = "  -2.50%  "

Set-TextValue $ws.Range("D22") "432.36"
$ws.Range("E22").Value = "  +3.88%  "

Set-TextValue $ws.Range("D23") "4.91"
$ws.Range("E23").Value = "  +11.08%  "

Set-TextValue $ws.Range("D24") "4.13"
$ws.Range("E24").Value = "  -4.50%  "

Set-TextValue $ws.Range("D25") "84.37"
$ws.Range("E25").Value = "  -1.57%  "

Set-TextValue $ws.Range("D26") "13.24"
$ws.Range("E26").Value = "  -0.99%  "

Set-TextValue $ws.Range("D27") "10.81"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("E28").Value = "  -0.46%  "

Set-TextValue $ws.Range("D29") "8.77"
$ws.Range("E29").Value = "  -3.80%  "

Set-TextValue $ws.Range("D30") "29.87"
$ws.Range("E30").Value = "  -1.75%  "

Set-TextValue $ws.Range("D31") "6.69"
$ws.Range("E31").Value = "  +2.05%  "

Set-TextValue $ws.Range("D32") "11.51"
$ws.Range("E32").Value = "  -2.14%  "

Set-TextValue $ws.Range("D33") "572.63"
$ws.Range("E33").Value = "  -5.70%  "

$ws.Range("E34").Value = "  -2.64%  "

Set-TextValue $ws.Range("D35") "58.48"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("E36").Value = "  +0.09%  "

Set-TextValue $ws.Range("D37") "0.143"
$ws.Range("E37").Value = "  -7.63%  "

Set-TextValue $ws.Range("D38") "3.51"
$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D39") "35.91"
$ws.Range("E39").Value = "  -3.61%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0762"
$ws.Range("E40").Value = "  -4.78%  "

Set-TextValue $ws.Range("D41") "0.370"
$ws.Range("E41").Value = "  -3.39%  "

$ws.Range("D42").Value = "3.121.48"
$ws.Range("E42").Value = "  -4.23%  "

Set-TextValue $ws.Range("D43") "0.997"
$ws.Range("E43").Value = "  -0.28%  "

Set-TextValue $ws.Range("D44") "2.84"
$ws.Range("E44").Value = "  -5.42%  "

Set-TextValue $ws.Range("D45") "3.27"
$ws.Range("E45").Value = "  -3.61%  "

Set-TextValue $ws.Range("D46") "0.0411"
$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("E47").Value = "  -2.89%  "

Set-TextValue $ws.Range("D48") "0.130"
$ws.Range("E48").Value = "  -1.76%  "

Set-TextValue $ws.Range("D49") "2.58"
$ws.Range("E49").Value = "  -3.87%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "8.28"
$ws.Range("E50").Value = "  -4.63%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "134.99"
$ws.Range("E51").Value = "  -2.44%  "
